$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A, shifting existing data right.
$ws.Range("A1:B1").EntireColumn.Insert()

# Match the header formatting used by the rest of row 1 (bold, centered, bordered).
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new header values.
$ws.Range("A1").Value = "Unnamed: 0.2"
$ws.Range("B1").Value = "Unnamed: 0.1"
